$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values on specific rows, per repulled data.
$ws.Range("F2").Value = -9
$ws.Range("F3").Value = 11
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 4
